$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 111, shifting rows 111:139 down to 112:140
$ws.Rows.Item(111).Insert()

# Populate the newly inserted row 111 with the new record
$ws.Cells.Item(111, 1).Value = 5
$ws.Cells.Item(111, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(111, 3).Value = "Maule"
$ws.Cells.Item(111, 4).Value = 44932
$ws.Cells.Item(111, 5).Value = 7
$ws.Cells.Item(111, 6).Value = 100112022
$ws.Cells.Item(111, 7).Value = "Arveja Verde"
$ws.Cells.Item(111, 8).Value = "Sin especificar"
$ws.Cells.Item(111, 9).Value = "Primera"
$ws.Cells.Item(111, 10).Value = 300
$ws.Cells.Item(111, 11).Value = 25000
$ws.Cells.Item(111, 12).Value = 25000
$ws.Cells.Item(111, 13).Value = 25000
$ws.Cells.Item(111, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(111, 15).Value = "Carahue"
$ws.Cells.Item(111, 16).Value = 1000
$ws.Cells.Item(111, 17).Value = 25
$ws.Cells.Item(111, 18).Value = "Hortaliza"
